$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the 3 rows belonging to the "CC / ROMAN RODRIGUEZ CAÑATE" worker.
#    This naturally shifts everything below (rows 19-25 -> 16-22, rows 30-31 -> 27-28)
#    along with merged cells, matching the target layout.
$ws.Rows("16:18").Delete()

# 2. Update header summary figures.
$ws.Range("E11").Value = 364000          # VALOR MORA
$ws.Range("C13").Value = 2               # Cant. Trabajadores
$ws.Range("F13").Value = 7               # Cant. Periodos

# 3. Update the worker document number for ANGELO SABIER VILLASMIL GARCIA.
$ws.Range("C16").Value = "20410188"

# 4. Re-order the "Periodo Mora" column for ENDER JOSE NAVAS FERREBUS rows
#    from descending (2309..2304) to ascending (2304..2309).
$ws.Range("E17").Value = "2304"
$ws.Range("E18").Value = "2305"
$ws.Range("E19").Value = "2306"
$ws.Range("E20").Value = "2307"
$ws.Range("E21").Value = "2308"
$ws.Range("E22").Value = "2309"
